$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the existing two "header" rows (rows 2 & 3) so that
# there is room for a third curfew-type header row. This shifts everything
# from row 4 onward down by one, which is exactly what we need since the
# bunk/staff rows used to start at row 5 and now need to start at row 6.
$ws.Rows.Item(4).Insert()

# Update the text of the (now three) merged curfew header rows.
$ws.Range("A2").Value = "Normal Curfew"
$ws.Range("A3").Value = "Night Off Curfew"
$ws.Range("A4").Value = "Day Off Curfew"

# Merge the newly inserted row the same way the other header rows are merged
# (drop the formatting copied down by Insert() first so the merge operation
# doesn't fork off a redundant duplicate of the centered style; re-apply the
# centered alignment afterwards so it matches rows 2 & 3 again).
$ws.Range("A4:G4").ClearFormats()
$ws.Range("A4:G4").Merge()
$ws.Range("A4:G4").HorizontalAlignment = -4108

# The "Staff Member 1 ID" cell used to be entered with a leading apostrophe
# (quote prefix) - reproduce that formatting.
$ws.Range("C6").Value = "'Staff Member 1 ID"

# Update the view's selection to match the saved state.
$ws.Range("H1:H1048576").Select()
